$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 ("Sheet1"): the two existing sample rows used the cable run
# "2C#2 / EXPRESS / 100+00 / 200+00" (row 3) and a mix for row 2. The pull
# sheet is being re-populated with 9 rows of real data: pulls 1-7 use
# "7C#14" and pulls 8-9 use "2C#2", all on the EXPRESS run from 100+00 to
# 200+00.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Clear the old sample rows first so the stale shared-string references to
# "2C#2", "100+00" and "200+00" are dropped before they get re-entered below
# (they need to be re-registered in the string table in the new row order).
$ws1.Range("A2:F3").ClearContents()

$cableSizes = @("7C#14","7C#14","7C#14","7C#14","7C#14","7C#14","7C#14","2C#2","2C#2")
for ($i = 0; $i -lt $cableSizes.Length; $i++) {
    $r = $i + 2
    $ws1.Cells.Item($r, 1).Value = $i + 1
    $ws1.Cells.Item($r, 2).Value = $cableSizes[$i]
    $ws1.Cells.Item($r, 3).Value = "EXPRESS"
    $ws1.Cells.Item($r, 4).Value = "100+00"
    $ws1.Cells.Item($r, 5).Value = "200+00"
}

[void]$ws1.Range("D11").Select()

# ---------------------------------------------------------------------------
# Sheet2 ("Sheet2"): add a new block of pull data (G5:K8) alongside the
# existing table, reusing the "2C#2 / EXPRESS / 100+00 / 200+00" run.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

$pullNumbers = @(3,4,5,6)
for ($i = 0; $i -lt $pullNumbers.Length; $i++) {
    $r = $i + 5
    $ws2.Cells.Item($r, 7).Value = $pullNumbers[$i]
    $ws2.Cells.Item($r, 8).Value = "2C#2"
    $ws2.Cells.Item($r, 9).Value = "EXPRESS"
    $ws2.Cells.Item($r, 10).Value = "100+00"
    $ws2.Cells.Item($r, 11).Value = "200+00"
}

[void]$ws2.Range("G5:K8").Select()
